# Update MeanLocPerEvent (column W) values for the re-processed rows.
# (Shared-string table de-duplication and workbook re-save are handled
# automatically by the engine when the file is written.)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("W2").Value  = 5.7803468208092479
$ws.Range("W3").Value  = 7.5757575757575761
$ws.Range("W4").Value  = 6.0606060606060606
$ws.Range("W5").Value  = 7.0175438596491224
$ws.Range("W6").Value  = 6.1728395061728394
$ws.Range("W7").Value  = 7.2992700729927007
$ws.Range("W9").Value  = 7.1174377224199281
$ws.Range("W10").Value = 7.7821011673151741
$ws.Range("W11").Value = 22.988505747126435
$ws.Range("W12").Value = 10.869565217391305
$ws.Range("W13").Value = 37.037037037037031
$ws.Range("W14").Value = 5.7636887608069163
$ws.Range("W15").Value = 6.2305295950155761
$ws.Range("W16").Value = 6.5573770491803272
$ws.Range("W17").Value = 5.9171597633136095
$ws.Range("W18").Value = 7.4626865671641784
$ws.Range("W20").Value = 5.7971014492753623
$ws.Range("W22").Value = 6.0790273556230998
$ws.Range("W23").Value = 6.1728395061728394
$ws.Range("W24").Value = 6.0790273556230998
$ws.Range("W25").Value = 6.968641114982578
